# Update cryptos list values per diff (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.256.74"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "2.296.32"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("D5").Value = "'315.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("D6").Value = "'102.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.15%  "
$ws.Range("D7").Value = "'0.622"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.56%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "'39.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("D11").Value = "'0.0902"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("E12").Value = "  +1.60%  "
$ws.Range("E13").Value = "  +0.68%  "
$ws.Range("D14").Value = "'0.960"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("D15").Value = "'15.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("D16").Value = "2.645.37"
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("D17").Value = "2.292.57"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").Value = "42.364.95"
$ws.Range("E18").Value = "  +1.18%  "
$ws.Range("D19").Value = "'7.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.79%  "
$ws.Range("D20").Value = "'0.0000105"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("D21").Value = "'73.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").Value = "'11.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +25.97%  "
$ws.Range("D23").Value = "'3.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.19%  "
$ws.Range("D24").Value = "'274.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.48%  "
$ws.Range("E25").Value = "  -2.17%  "
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").Value = "'10.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("E28").Value = "  +3.50%  "
$ws.Range("D29").Value = "'22.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("D30").Value = "'37.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.01%  "
$ws.Range("D31").Value = "'165.65"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").Value = "'0.0874"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("D33").Value = "'5.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.68%  "
$ws.Range("D34").Value = "'0.133"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.40%  "
$ws.Range("E35").Value = "  -9.70%  "
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").Value = "'4.57"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("D38").Value = "'0.0363"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.28%  "
$ws.Range("D39").Value = "'3.72"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.72%  "
$ws.Range("D40").Value = "'2.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.05%  "
$ws.Range("D41").Value = "'1.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.51%  "
$ws.Range("D42").Value = "'70.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.10%  "
$ws.Range("D43").Value = "'94.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.97%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "'0.225"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("D46").Value = "'12.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.69%  "
$ws.Range("D47").Value = "'80.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.40%  "
$ws.Range("D48").Value = "'112.91"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.76%  "
$ws.Range("D49").Value = "'8.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("D51").Value = "1.589.43"
$ws.Range("E51").Value = "  +1.98%  "